$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 (Serie = 01-01-2021) with revised values
$ws.Range("B74").Value = -178
$ws.Range("C74").Value = -87
$ws.Range("D74").Value = 32
$ws.Range("E74").Value = -119
$ws.Range("F74").Value = -12
$ws.Range("G74").Value = 58
$ws.Range("H74").Value = -69
$ws.Range("I74").Value = -18
$ws.Range("J74").Value = -18
$ws.Range("K74").Value = 149
$ws.Range("L74").Value = -47
$ws.Range("M74").Value = -31
$ws.Range("N74").Value = 227
$ws.Range("O74").Value = -211
$ws.Range("P74").Value = 81
$ws.Range("Q74").Value = -259
$ws.Range("R74").Value = -144
$ws.Range("S74").Value = -103
$ws.Range("T74").Value = -41
$ws.Range("U74").Value = 12
$ws.Range("V74").Value = 12
$ws.Range("W74").Value = 40
$ws.Range("X74").Value = 44
$ws.Range("Y74").Value = -4
$ws.Range("Z74").Value = -167

# Add a new row 75 for the next quarter (01-04-2021)
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = -686
$ws.Range("C75").Value = -95
$ws.Range("D75").Value = -31
$ws.Range("E75").Value = -64
$ws.Range("F75").Value = -663
$ws.Range("G75").Value = -35
$ws.Range("H75").Value = -628
$ws.Range("I75").Value = 136
$ws.Range("J75").Value = 136
$ws.Range("K75").Value = -196
$ws.Range("L75").Value = -106
$ws.Range("M75").Value = -20
$ws.Range("N75").Value = -69
$ws.Range("O75").Value = 132
$ws.Range("P75").Value = -93
$ws.Range("Q75").Value = -592
$ws.Range("R75").Value = 119
$ws.Range("S75").Value = 88
$ws.Range("T75").Value = 31
$ws.Range("U75").Value = 41
$ws.Range("V75").Value = 41
$ws.Range("W75").Value = -560
$ws.Range("X75").Value = -520
$ws.Range("Y75").Value = -40
$ws.Range("Z75").Value = -193
